{"js": "// Rename three entries in the \"Forms\" table: fix the Korean technique\n// names \"Kicho E Jang\" -> \"Kicho I Jang\", \"Kicho Sam Jan\" -> \"Kicho Sam Jang\",\n// and \"Palgue E Jang\" -> \"Palgue I Jang\".\nconst replacements = [\n  [\"Kicho E Jang\", \"Kicho I Jang\"],\n  [\"Kicho Sam Jan\", \"Kicho Sam Jang\"],\n  [\"Palgue E Jang\", \"Palgue I Jang\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Rename three entries in the \"Forms\" table: fix the Korean technique\n# names \"Kicho E Jang\" -> \"Kicho I Jang\", \"Kicho Sam Jan\" -> \"Kicho Sam Jang\",\n# and \"Palgue E Jang\" -> \"Palgue I Jang\".\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Kicho E Jang\", \"Kicho I Jang\"),\n    @(\"Kicho Sam Jan\", \"Kicho Sam Jang\"),\n    @(\"Palgue E Jang\", \"Palgue I Jang\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
